$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Cfh", "Sell", "ECs", 3, 1, 1.573304666666667, 4.719914, 0.02227938362376841, 0.0222793836237684, 1, 0.3333333333333333, 74.89537566666667, 224.686127, 0.6660623326691122, 0.6660623326691123, 117.8332440481198, 1060.499196433078, 0.0148394582268772, 0.0148394582268772),
    @("ECs", "Cfh", "Sell", "M2", 3, 1, 1.573304666666667, 4.719914, 0.02227938362376841, 0.0222793836237684, 3, 1, 37.54961933333333, 112.648858, 0.3339376673308877, 0.3339376673308878, 59.07699132869023, 531.692921958212, 0.007439925396891202, 0.007439925396891202),
    @("FAPs", "Cfh", "Sell", "ECs", 3, 1, 36.796554, 110.389662, 0.5210717033819111, 0.521071703381911, 1, 0.3333333333333333, 74.89537566666667, 224.686127, 0.6660623326691122, 0.6660623326691123, 2755.891735068786, 24803.02561561907, 0.3470662342424234, 0.3470662342424234),
    @("FAPs", "Cfh", "Sell", "M2", 3, 1, 36.796554, 110.389662, 0.5210717033819111, 0.521071703381911, 3, 1, 37.54961933333333, 112.648858, 0.3339376673308877, 0.3339376673308878, 1381.696595478444, 12435.269359306, 0.1740054691394876, 0.1740054691394876),
    @("M2", "Cfh", "Sell", "ECs", 3, 1, 17.68148033333333, 53.04444099999999, 0.2503853778156443, 0.2503853778156443, 1, 0.3333333333333333, 74.89537566666667, 224.686127, 0.6660623326691122, 0.6660623326691123, 1324.261111907778, 11918.35000717001, 0.166772268814125, 0.1667722688141251),
    @("M2", "Cfh", "Sell", "M2", 3, 1, 17.68148033333333, 53.04444099999999, 0.2503853778156443, 0.2503853778156443, 3, 1, 37.54961933333333, 112.648858, 0.3339376673308877, 0.3339376673308878, 663.9328557664862, 5975.395701898377, 0.08361310900151928, 0.08361310900151929),
    @("sCs", "Cfh", "Sell", "ECs", 3, 1, 14.56572533333333, 43.697176, 0.2062635351786761, 0.2062635351786761, 1, 0.3333333333333333, 74.89537566666667, 224.686127, 0.6660623326691122, 0.6660623326691123, 1090.905470697483, 9818.149236277352, 0.1373843713856865, 0.1373843713856865),
    @("sCs", "Cfh", "Sell", "M2", 3, 1, 14.56572533333333, 43.697176, 0.2062635351786761, 0.2062635351786761, 3, 1, 37.54961933333333, 112.648858, 0.3339376673308877, 0.3339376673308878, 546.9374415805564, 4922.436974225007, 0.0688791637929896, 0.0688791637929896)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}